$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("leadlag")

$row = 30

# Trade #
$ws.Cells.Item($row, 1).Value = 34

# Date / Time - force text format so "2026-02-16" / "21:28:34" aren't
# auto-converted into date/time serial numbers, matching the rest of the sheet,
# then restore the "Normal" style so no stray number format sticks to the cell
$ws.Cells.Item($row, 2).NumberFormat = "@"
$ws.Cells.Item($row, 2).Value = "2026-02-16"
$ws.Cells.Item($row, 2).Style = "Normal"

$ws.Cells.Item($row, 3).NumberFormat = "@"
$ws.Cells.Item($row, 3).Value = "21:28:34"
$ws.Cells.Item($row, 3).Style = "Normal"

# Strategy / Side
$ws.Cells.Item($row, 4).Value = "leadlag"
$ws.Cells.Item($row, 5).Value = "DOWN"

# Entry Price
$ws.Cells.Item($row, 6).Value = 68810.35000000001

# Exit Price - left blank (trade still open)
$ws.Cells.Item($row, 7).Value = ""

# Status
$ws.Cells.Item($row, 8).Value = "OPEN"

# P&L % / P&L $
$ws.Cells.Item($row, 9).Value = 0
$ws.Cells.Item($row, 10).Value = 0

# Confidence
$ws.Cells.Item($row, 11).Value = 0.75

# Entry Reason
$ws.Cells.Item($row, 12).Value = "Binance leading with -0.125% move"

# Exit Reason - left blank (trade still open)
$ws.Cells.Item($row, 13).Value = ""

# Duration (min)
$ws.Cells.Item($row, 14).Value = 0
